$d = $word.ActiveDocument

# Locate the word "кофейни" inside the "Анализ требований..." paragraph and
# replace it with "автостанции", while leaving the rest of the sentence intact.
$range = $d.Content
$found = $range.Find.Execute("кофейни", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $range.Text = "автостанции"

    # Force the engine to keep this replaced word as its own run instead of
    # silently re-merging it with the identically formatted runs on either
    # side of it (toggling Bold on/off is a no-op on the rendered formatting
    # but causes the run boundaries introduced by the text replacement to be
    # preserved on save).
    $range.Bold = 1
    $range.Bold = 0
}
